$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = -1
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = 5
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = -2
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = 4
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = -3
$ws.Range("F24").Value = 2
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 5
$ws.Range("F27").Value = -1
$ws.Range("F28").Value = -1
$ws.Range("F32").Value = -1
